# Update team-specific time data values in Sheet1 (Auburn_B)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1871657754010695
$ws.Range("C2").Value = 0.5721925133689839
$ws.Range("J2").Value = 0.01336898395721925
$ws.Range("O2").Value = 0.00267379679144385
$ws.Range("P2").Value = 0.1417112299465241
$ws.Range("S2").Value = 0.08288770053475936

# Row 3
$ws.Range("B3").Value = 0.00909090909090909
$ws.Range("C3").Value = 0.02727272727272727
$ws.Range("J3").Value = 0.03181818181818181
$ws.Range("P3").Value = 0.7636363636363637
$ws.Range("S3").Value = 0.1681818181818182

# Row 4
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.6595744680851063
$ws.Range("S4").Value = 0.3191489361702128

# Row 6
$ws.Range("B6").Value = 0.06880733944954129
$ws.Range("D6").Value = 0.01376146788990826
$ws.Range("F6").Value = 0.1238532110091743
$ws.Range("J6").Value = 0.1743119266055046
$ws.Range("O6").Value = 0.05045871559633028
$ws.Range("Q6").Value = 0.1559633027522936
$ws.Range("R6").Value = 0.06422018348623854
$ws.Range("S6").Value = 0.3486238532110092

# Row 7
$ws.Range("B7").Value = 0.130952380952381
$ws.Range("D7").Value = 0.01587301587301587
$ws.Range("F7").Value = 0.03571428571428571
$ws.Range("J7").Value = 0.123015873015873
$ws.Range("O7").Value = 0.007936507936507936
$ws.Range("Q7").Value = 0.130952380952381
$ws.Range("R7").Value = 0.06349206349206349
$ws.Range("S7").Value = 0.492063492063492

# Row 8
$ws.Range("B8").Value = 0.096579476861167
$ws.Range("D8").Value = 0.01810865191146881
$ws.Range("F8").Value = 0.06237424547283702
$ws.Range("J8").Value = 0.1207243460764588
$ws.Range("O8").Value = 0.01810865191146881
$ws.Range("Q8").Value = 0.1428571428571428
$ws.Range("R8").Value = 0.07847082494969819
$ws.Range("S8").Value = 0.4627766599597585

# Row 9
$ws.Range("B9").Value = 0.1173708920187793
$ws.Range("D9").Value = 0.01408450704225352
$ws.Range("F9").Value = 0.0892018779342723
$ws.Range("J9").Value = 0.1643192488262911
$ws.Range("O9").Value = 0.004694835680751174
$ws.Range("Q9").Value = 0.1408450704225352
$ws.Range("R9").Value = 0.08450704225352113
$ws.Range("S9").Value = 0.3849765258215962

# Row 10
$ws.Range("B10").Value = 0.1414868105515588
$ws.Range("D10").Value = 0.0231814548361311
$ws.Range("F10").Value = 0.05995203836930456
$ws.Range("J10").Value = 0.158273381294964
$ws.Range("O10").Value = 0.01998401278976818
$ws.Range("Q10").Value = 0.1814548361310951
$ws.Range("R10").Value = 0.06554756195043965
$ws.Range("S10").Value = 0.3501199040767386

# Row 11
$ws.Range("G11").Value = 0.1282798833819242
$ws.Range("J11").Value = 0.07580174927113703
$ws.Range("K11").Value = 0.1720116618075802
$ws.Range("L11").Value = 0.6064139941690962
$ws.Range("S11").Value = 0.01749271137026239

# Row 12
$ws.Range("G12").Value = 0.7863636363636364
$ws.Range("J12").Value = 0.15
$ws.Range("K12").Value = 0.004545454545454545
$ws.Range("L12").Value = 0.03181818181818181
$ws.Range("S12").Value = 0.02727272727272727

# Row 13
$ws.Range("G13").Value = 0.676923076923077
$ws.Range("J13").Value = 0.2615384615384616
$ws.Range("S13").Value = 0.06153846153846154

# Row 15
$ws.Range("F15").Value = 0.02631578947368421
$ws.Range("H15").Value = 0.1885964912280702
$ws.Range("I15").Value = 0.07456140350877193
$ws.Range("J15").Value = 0.2456140350877193
$ws.Range("K15").Value = 0.09649122807017543
$ws.Range("M15").Value = 0.03070175438596491
$ws.Range("O15").Value = 0.09210526315789473
$ws.Range("S15").Value = 0.2456140350877193

# Row 16
$ws.Range("F16").Value = 0.025
$ws.Range("H16").Value = 0.2333333333333333
$ws.Range("I16").Value = 0.05416666666666667
$ws.Range("J16").Value = 0.3041666666666666
$ws.Range("K16").Value = 0.1416666666666667
$ws.Range("M16").Value = 0.03333333333333333
$ws.Range("O16").Value = 0.07083333333333333
$ws.Range("S16").Value = 0.1375

# Row 17
$ws.Range("F17").Value = 0.01269035532994924
$ws.Range("H17").Value = 0.1649746192893401
$ws.Range("I17").Value = 0.1040609137055838
$ws.Range("J17").Value = 0.4492385786802031
$ws.Range("K17").Value = 0.08375634517766498
$ws.Range("M17").Value = 0.03045685279187817
$ws.Range("N17").Value = 0.002538071065989848
$ws.Range("O17").Value = 0.0583756345177665
$ws.Range("S17").Value = 0.09390862944162437

# Row 18
$ws.Range("F18").Value = 0.01744186046511628
$ws.Range("H18").Value = 0.1744186046511628
$ws.Range("I18").Value = 0.09883720930232558
$ws.Range("J18").Value = 0.4011627906976744
$ws.Range("K18").Value = 0.09302325581395349
$ws.Range("M18").Value = 0.005813953488372093
$ws.Range("N18").Value = 0.005813953488372093
$ws.Range("O18").Value = 0.06976744186046512
$ws.Range("S18").Value = 0.1337209302325581

# Row 19
$ws.Range("F19").Value = 0.01251840942562592
$ws.Range("H19").Value = 0.2187039764359352
$ws.Range("I19").Value = 0.0898379970544919
$ws.Range("J19").Value = 0.3321060382916053
$ws.Range("K19").Value = 0.1303387334315169
$ws.Range("M19").Value = 0.02798232695139912
$ws.Range("N19").Value = 0.002209131075110457
$ws.Range("O19").Value = 0.05522827687776141
$ws.Range("S19").Value = 0.1310751104565538

$wb.Save()